# The workbook stores a demo login on Sheet1:
#   A1=admin            B1=demo
#   A2=validuname (*)   B2=validpwd
# (*) A2 used to carry a mailto: hyperlink to santhiparambalam@gmail.com.
#
# This edit swaps in the real test credentials, drops the obsolete
# hyperlink (while keeping the cell's existing "Hyperlink" look), and
# moves the active selection off of B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the outdated mailto hyperlink that lived on A2.
$ws.Range("A2").Hyperlinks.Delete()

# Write the new username value into A2 via copy/paste-values so the
# cell's existing formatting (the old Hyperlink style) is left intact
# instead of being replaced by a fresh direct-format copy.
$ws.Range("D1").Value = "santhi.asusvivobook@gmail.com"
$ws.Range("D1").Copy()
$ws.Range("A2").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("D1").ClearContents()
$excel.CutCopyMode = $false

# Update the password value.
$ws.Range("B2").Value = "santhi_asusvivobook"

# Move the selection from B2 to F2.
$ws.Range("F2").Select()
